# Remove the "p_timeStep_h" parameter row from the params sheet and replace
# the remaining row's content with the "p_undergroundTemperature_degC" row,
# so the timestep is no longer sourced from the local workbook (it is now
# driven by the experiment instead).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite row 2 (previously p_timeStep_h) with the former row 3 data
# (p_undergroundTemperature_degC), dropping the now-unused "comment" value.
$ws.Range("B2").Value = "p_undergroundTemperature_degC"
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = "degC"
$ws.Range("E2").ClearContents() | Out-Null

# Remove the old row 3 entirely, shifting the used range up.
$ws.Rows("3:3").Delete() | Out-Null

# Match the author's last-selected cell.
$ws.Range("B3").Select() | Out-Null
